$wb = $excel.ActiveWorkbook

# --- Update p-values ("rerun with v12 hist") ---
# Cells are addressed by their *original* tab name (sheetId/r:id-stable)
# so the later tab-name swap below cannot affect which sheet gets which value.
$ws = $wb.Worksheets.Item("Low-grade glioma")
$ws.Range("C2").Value = 0.15230869731919
$ws.Range("C4").Value = 0.7931779162466
$ws.Range("C7").Value = 0.869185449350616
$ws.Range("C8").Value = 0.62685469729387
$ws.Range("C9").Value = 0.76454731048697

$ws = $wb.Worksheets.Item("Medulloblastoma")
$ws.Range("C7").Value = 0.0231617439105454

$ws = $wb.Worksheets.Item("Mixed neuronal-glial tumor")
$ws.Range("C2").Value = 0.511414611707042
$ws.Range("C3").Value = 0.0191251203686029
$ws.Range("C4").Value = 0.830125478641495
$ws.Range("C5").Value = 0.608353352840895
$ws.Range("C6").Value = 0.532508089844862
$ws.Range("C7").Value = 0.190274200316843
$ws.Range("C8").Value = 0.193583517235016
$ws.Range("C9").Value = 0.367827062567279

$ws = $wb.Worksheets.Item("Ependymoma")
$ws.Range("C2").Value = 0.0485949954198866
$ws.Range("C4").Value = 0.357593243242194
$ws.Range("C5").Value = 0.212651737225106
$ws.Range("C6").Value = 0.384467290382669
$ws.Range("C7").Value = 0.212297500884979
$ws.Range("C8").Value = 0.0335746893173699
$ws.Range("C9").Value = 0.768522422425236

$ws = $wb.Worksheets.Item("Other high-grade glioma")
$ws.Range("C2").Value = 0.197737271253414
$ws.Range("C3").Value = 0.0860576559380436
$ws.Range("C4").Value = 0.796815847884649
$ws.Range("C5").Value = 0.844003973383851
$ws.Range("C6").Value = 0.999999999999984
$ws.Range("C7").Value = 0.442905776000836
$ws.Range("C8").Value = 1
$ws.Range("C9").Value = 0.299113352181795

$ws = $wb.Worksheets.Item("Craniopharyngioma")
$ws.Range("C4").Value = 0.644381223328592

$ws = $wb.Worksheets.Item("Mesenchymal tumor")
$ws.Range("C2").Value = 0.0791739604170237
$ws.Range("C3").Value = 1
$ws.Range("C4").Value = 0.746180752659459
$ws.Range("C5").Value = 1
$ws.Range("C6").Value = 0.907161981862534
$ws.Range("C7").Value = 0.344182227681443
$ws.Range("C8").Value = 0.613245356793744
$ws.Range("C9").Value = 0.00457730524302

$ws = $wb.Worksheets.Item("DIPG or DMG")
$ws.Range("C2").Value = 0.681940144478844
$ws.Range("C3").Value = 0.284829721362229
$ws.Range("C4").Value = 1
$ws.Range("C5").Value = 0.508771929824559
$ws.Range("C6").Value = 0.999999999999994
$ws.Range("C7").Value = 0.733266733266733
$ws.Range("C8").Value = 0.0759240759240759
$ws.Range("C9").Value = 0.635711125445117

$ws = $wb.Worksheets.Item("Neurofibroma plexiform")
$ws.Range("C2").Value = 0.297702297702298

$ws = $wb.Worksheets.Item("Non-neoplastic tumor")
$ws.Range("C2").Value = 0.422874973996255
$ws.Range("C3").Value = 0.659246931558144
$ws.Range("C4").Value = 0.564584980237159
$ws.Range("C5").Value = 0.564584980237159
$ws.Range("C6").Value = 0.564584980237159
$ws.Range("C7").Value = 0.43956043956044
$ws.Range("C9").Value = 0.278072325213818

$ws = $wb.Worksheets.Item("Schwannoma")
$ws.Range("C2").Value = 0.33006993006993

$ws = $wb.Worksheets.Item("Other tumor")
$ws.Range("C2").Value = 0.387996466662305
$ws.Range("C3").Value = 0.389108813714585
$ws.Range("C4").Value = 0.275804815808414
$ws.Range("C5").Value = 0.999999999999982
$ws.Range("C6").Value = 0.762971929594966
$ws.Range("C7").Value = 0.721120984278879
$ws.Range("C8").Value = 0.626081012808702
$ws.Range("C9").Value = 0.403692093453513

# --- Swap the "Mesenchymal tumor" / "DIPG or DMG" tab names ---
# (sheetId/r:id - i.e. the underlying data - stay put; only the visible
# tab names trade places, per the workbook.xml <sheets> diff.)
$wsA = $wb.Worksheets.Item("Mesenchymal tumor")
$wsB = $wb.Worksheets.Item("DIPG or DMG")
$wsA.Name = "__tmp_swap__"
$wsB.Name = "Mesenchymal tumor"
$wsA.Name = "DIPG or DMG"
